$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.729.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +7.36%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.503.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +8.51%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.39%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "191.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +14.93%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "556.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.05%  "
# Row 7
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.495.22"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.24%  "
# Row 8
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.614"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.76%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.21%  "
# Row 10
$ws.Range("E10").Value = "  +8.44%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.93%  "
# Row 12
$ws.Range("E12").Value = "  +16.68%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.75%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.01%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.060.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.51%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.502.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.72%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.945.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.77%  "
# Row 18
$ws.Range("E18").Value = "  +6.24%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.70%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.23%  "
# Row 21
$ws.Range("E21").Value = "  +8.27%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "406.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.03%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.40%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.31%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.47%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.77%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.60%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.35%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.81%  "
# Row 30
$ws.Range("E30").Value = "  +7.71%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "692.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.55%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.13%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.62%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.10%  "
# Row 35
$ws.Range("E35").Value = "  +8.56%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "60.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.67%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.89%  "
# Row 38
$ws.Range("E38").Value = "  +24.88%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.404"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.60%  "
# Row 40
$ws.Range("E40").Value = "  -0.16%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +26.79%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.134"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.68%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +18.50%  "
# Row 44
$ws.Range("E44").Value = "  +0.17%  "
# Row 45
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.25%  "
# Row 46
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.055.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.83%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0423"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.82%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.33%  "
# Row 49
$ws.Range("E49").Value = "  +9.64%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +16.28%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.130"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.18%  "
